# Upgrade to ClosedXML 0.101 (#321)
# "Fix test that started to fail after upgrading calc engine in ClosedXML;
#  make the test culture-independent"
#
# The B3 cell on sheet "Лист1" used a volatile, locale-sensitive formula
# (TODAY()/TEXT(...,"dd.MM.yyyy")) that produced different results depending
# on the day it was evaluated and the culture in use. Replace it with a
# deterministic, culture-independent formula that hard-codes the expected
# date string, matching what the upgraded calc engine now expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("B3").Formula = '=CONCATENATE("Begin at ","19.01.2023")'
